# git-mingling.xlsx refresh: add new rows 14-34 with additional git commands
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# helper: apply the existing "column C" style (s="1", fontId 3 = 等线/family3/charset134)
# by copying format from an already-styled cell (C2) - this reuses the existing
# style index instead of fabricating a new (slightly different) font.
function Apply-ColCStyle($range) {
    $ws.Range("C2").Copy() | Out-Null
    $range.PasteSpecial(-4122) | Out-Null
}

$excel.CutCopyMode = 0

# ---- Row 14 ----
$ws.Range("B14").Value = "git rev-parse --git-dir"
$ws.Range("C14").Value = "显示.git目录所在的位置"
Apply-ColCStyle $ws.Range("C14")

# ---- Row 15 ----
$ws.Range("B15").Value = "git rev-parse --show-toplevel"
$ws.Range("C15").Value = "显示工作区的根目录"
Apply-ColCStyle $ws.Range("C15")

# ---- Row 16 ----
$ws.Range("B16").Value = "git rev-parse --show-prefix"
$ws.Range("C16").Value = "相对于工作区的相对目录"
Apply-ColCStyle $ws.Range("C16")

# ---- Row 17 ----
$ws.Range("B17").Value = "git config -e --global"
$ws.Range("C17").Value = "编辑全局配置文件"
Apply-ColCStyle $ws.Range("C17")

# ---- Row 18 ----
$ws.Range("B18").Value = "git config -e --system"
$ws.Range("C18").Value = "编辑系统配置文件"
Apply-ColCStyle $ws.Range("C18")

# ---- Row 19 ----
$ws.Range("B19").Value = "git config --unset --global user.name"
$ws.Range("C19").Value = "删除全局文件中的user.name"
Apply-ColCStyle $ws.Range("C19")
$ws.Range("C19").Characters(9, 9).Font.Name = "等线"

# ---- Row 20 (re-uses the exact same C text/rich-run as row 19) ----
$ws.Range("B20").Value = "git config --unset --global user.email"
$ws.Range("C20").Value = "删除全局文件中的user.name"
Apply-ColCStyle $ws.Range("C20")
$ws.Range("C20").Characters(9, 9).Font.Name = "等线"

# ---- Row 21 ----
$ws.Range("B21").Value = "git config user.name"
$ws.Range("C21").Value = "查看全局文件中的用户名"
Apply-ColCStyle $ws.Range("C21")

# ---- Row 22 ----
$ws.Range("B22").Value = "git config user.email"
$ws.Range("C22").Value = "查看全局文件中的用户邮件"
Apply-ColCStyle $ws.Range("C22")

# ---- Row 23 ----
$ws.Range("B23").Value = 'git commit --allow-empty -m "msg"'
$ws.Range("C23").Value = "允许空白提交（无任何文件修改）"
Apply-ColCStyle $ws.Range("C23")

# ---- Row 24 ----
$ws.Range("B24").Value = "git log --pretty=fuller"
$ws.Range("C24").Value = "日志全显"
Apply-ColCStyle $ws.Range("C24")

# ---- Row 25 (quote-prefixed cell: leading apostrophe => quotePrefix style) ----
$ws.Range("B25").Value = "git commit --amend --allow-empty --reset-author"
$ws.Range("C25").Value = "'--amend对刚刚提交进行修补"
$ws.Range("C25").Characters(8, 9).Font.Name = "等线"

# ---- Row 26 (quote-prefixed, column-C font) ----
$ws.Range("C26").Value = "'--allow-empty使得空表提交被允许"
Apply-ColCStyle $ws.Range("C26")

# ---- Row 27 (starts with a curly quote U+2019, NOT an apostrophe -> no quote-prefix) ----
$ws.Range("C27").Value = "’--reset-author将Author的ID同步修改"
Apply-ColCStyle $ws.Range("C27")
$ws.Range("C27").Characters(10, 6).Font.Name = "等线"
$ws.Range("C27").Characters(17, 2).Font.Name = "等线"

# ---- Row 28 ----
$ws.Range("B28").Value = "git push origin master"
$ws.Range("C28").Value = "master为分支，提交到github中"
Apply-ColCStyle $ws.Range("C28")

# ---- Row 29 ----
$ws.Range("B29").Value = "git log --state"
$ws.Range("C29").Value = "查看提交日志，--state可以看到每次提交的文件变更统计"
Apply-ColCStyle $ws.Range("C29")
$ws.Range("C29").Characters(8, 7).Font.Name = "等线"

# ---- Row 30 ----
$ws.Range("B30").Value = "git diff"
$ws.Range("C30").Value = "查看修改后的文件与版本库中的文件差异"
Apply-ColCStyle $ws.Range("C30")

# ---- Row 31 ----
$ws.Range("B31").Value = "git status -s"
$ws.Range("C31").Value = "查看文件信息，-s简化信息量"
Apply-ColCStyle $ws.Range("C31")
$ws.Range("C31").Characters(8, 2).Font.Name = "等线"

# ---- Row 32 ----
$ws.Range("B32").Value = "git checkout -- filename"
$ws.Range("C32").Value = "撤销工作区尚未提交的的修改"
Apply-ColCStyle $ws.Range("C32")

# ---- Row 33 ----
$ws.Range("B33").Value = "git ls-tree -l HEAD"
$ws.Range("C33").Value = "查看暂存区及HEAD中的目录树"
Apply-ColCStyle $ws.Range("C33")
$ws.Range("C33").Characters(7, 4).Font.Name = "等线"

# ---- Row 34 ----
$ws.Range("B34").Value = "git clean -fd"
$ws.Range("C34").Value = "清除当前工作区中没有加入版本库的文件和目录"
Apply-ColCStyle $ws.Range("C34")

$excel.CutCopyMode = 0

# ---- sheet view: scroll + selection ----
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 2
$ws.Range("B35").Select()

# ---- page setup (adds pageSetup element, orientation=portrait) ----
$ws.PageSetup.Orientation = 1
